# [Kadastro App] Yeni kayit eklendi: 2961
# Adds the new record (Kayit No 2961) to the "Kayitlar" master sheet and to
# the "Erdemli" filtered sheet (the record's Birim/region).

$wb = $excel.ActiveWorkbook

function Add-KayitRow {
    param($ws, $row)

    $rng = $ws.Range("A" + $row + ":F" + $row)
    # Force the whole new row to be stored as text, matching the rest of the
    # sheet (every cell in these sheets is text, even numeric-looking ones,
    # with an ignoredErrors/numberStoredAsText marker covering the range).
    $rng.NumberFormat = "@"

    $ws.Range("A" + $row).Value = "2961"
    $ws.Range("B" + $row).Value = "2025-09-09"
    $ws.Range("C" + $row).Value = "Erdemli"
    $ws.Range("D" + $row).Value = ""
    $ws.Range("E" + $row).Value = "3B"
    $ws.Range("F" + $row).Value = "ÖZKAN AKBAŞ (Mühendis), SEVİL SARAÇER (Tekniker)"

    # Drop the temporary "@" number format again so the new cells end up
    # styled exactly like the rest of the (unstyled) sheet.
    $rng.ClearFormats()
}

# "Kayitlar" is the master log sheet: new record becomes row 31.
$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
Add-KayitRow $wsKayitlar 31

# "Erdemli" is the per-region sheet filtered to Birim = Erdemli: new record
# becomes row 30 there (one fewer row than the master sheet, since the
# master sheet also contains an Anamur record not present here).
$wsErdemli = $wb.Worksheets.Item("Erdemli")
Add-KayitRow $wsErdemli 30
